# Update EPC data (combined, every-3-days)
# Appends a new row (row 4) of data to both the "Finance" and "Non-Finance"
# worksheets, mirroring the structure already present in row 3 / row 2.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Finance" ----
$wsFinance = $wb.Worksheets.Item("Finance")

# Date column - must stay plain text (not get auto-converted to a date
# serial number), so we temporarily force a text number format, set the
# value, then clear the formatting back to the default/general style.
$cellA = $wsFinance.Range("A4")
$cellA.NumberFormat = "@"
$cellA.Value = "2025-10-22"
$cellA.ClearFormats()

$wsFinance.Range("B4").Value = 5.68
$wsFinance.Range("C4").Value = 9.32
$wsFinance.Range("D4").Value = 7.66
$wsFinance.Range("E4").Value = 5.03
$wsFinance.Range("F4").Value = 26.84
$wsFinance.Range("G4").Value = 2.73
$wsFinance.Range("H4").Value = 5.36
$wsFinance.Range("I4").Value = 3.5
$wsFinance.Range("J4").Value = 2.73
$wsFinance.Range("K4").Value = 12.58
$wsFinance.Range("L4").Value = 5.83
$wsFinance.Range("M4").Value = 5.68

# ---- Sheet 2: "Non-Finance" ----
$wsNonFinance = $wb.Worksheets.Item("Non-Finance")

$cellA2 = $wsNonFinance.Range("A4")
$cellA2.NumberFormat = "@"
$cellA2.Value = "2025-10-22"
$cellA2.ClearFormats()
